# Add a new "2022-Q4" sheet (positioned right after "总计", before "2022-Q3")
# and update the "总计" summary sheet with the new quarter's totals.
# The existing "2022-Q3" / "2022-Q2" / "2022-Q1" sheets are left untouched —
# they simply shift right by one tab position because of the new sheet.

$wb = $excel.ActiveWorkbook

$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ3 = $wb.Worksheets.Item(2)

# --- Create the new "2022-Q4" sheet by copying "2022-Q3" (keeps headers/styles) ---
$sheetQ3.Copy($sheetQ3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2: 160910 大成创新成长混合（LOF）
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "160910"
$q4.Cells.Item(2,3).Value = "大成创新成长混合（LOF）"
$q4.Cells.Item(2,4).Value = "14.01"
$q4.Cells.Item(2,5).Value = "78.14"
$q4.Cells.Item(2,6).Value = "6.50"
$q4.Cells.Item(2,7).Value = "0.9106"
$q4.Cells.Item(2,8).Value = 2

# Row 3: 010826 大成产业趋势混合A
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "'010826"
$q4.Cells.Item(3,3).Value = "大成产业趋势混合A"
$q4.Cells.Item(3,4).Value = "11.37"
$q4.Cells.Item(3,5).Value = "93.99"
$q4.Cells.Item(3,6).Value = "3.47"
$q4.Cells.Item(3,7).Value = "0.3945"
$q4.Cells.Item(3,8).Value = 10

# Row 4: 010827 大成产业趋势混合C (new row, copy column-A style from row 2)
$q4.Cells.Item(2,1).Copy()
$q4.Cells.Item(4,1).PasteSpecial(-4122)
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "'010827"
$q4.Cells.Item(4,3).Value = "大成产业趋势混合C"
$q4.Cells.Item(4,4).Value = "3.42"
$q4.Cells.Item(4,5).Value = "93.99"
$q4.Cells.Item(4,6).Value = "3.47"
$q4.Cells.Item(4,7).Value = "0.1187"
$q4.Cells.Item(4,8).Value = 10

# Row 5: 233009 大摩多因子精选策略混合 (new row)
$q4.Cells.Item(2,1).Copy()
$q4.Cells.Item(5,1).PasteSpecial(-4122)
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "233009"
$q4.Cells.Item(5,3).Value = "大摩多因子精选策略混合"
$q4.Cells.Item(5,4).Value = "6.42"
$q4.Cells.Item(5,5).Value = "91.11"
$q4.Cells.Item(5,6).Value = "0.98"
$q4.Cells.Item(5,7).Value = "0.0629"
$q4.Cells.Item(5,8).Value = 5

# Row 6: 015707 安信新能源主题股票A (new row)
$q4.Cells.Item(2,1).Copy()
$q4.Cells.Item(6,1).PasteSpecial(-4122)
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "'015707"
$q4.Cells.Item(6,3).Value = "安信新能源主题股票A"
$q4.Cells.Item(6,4).Value = "0.18"
$q4.Cells.Item(6,5).Value = "87.28"
$q4.Cells.Item(6,6).Value = "3.47"
$q4.Cells.Item(6,7).Value = "0.0062"
$q4.Cells.Item(6,8).Value = 10

# Row 7: 015708 安信新能源主题股票C (was row 4 on "2022-Q3"; updated values)
$q4.Cells.Item(2,1).Copy()
$q4.Cells.Item(7,1).PasteSpecial(-4122)
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "'015708"
$q4.Cells.Item(7,3).Value = "安信新能源主题股票C"
$q4.Cells.Item(7,4).Value = "0.13"
$q4.Cells.Item(7,5).Value = "87.28"
$q4.Cells.Item(7,6).Value = "3.47"
$q4.Cells.Item(7,7).Value = "0.0045"
$q4.Cells.Item(7,8).Value = 10

# --- Update the "总计" sheet with the new Q4 summary row, shifting the others down ---
$sheetTotal.Cells.Item(2,1).Value = 0
$sheetTotal.Cells.Item(2,2).Value = "2022-Q4"
$sheetTotal.Cells.Item(2,3).Value = 6
$sheetTotal.Cells.Item(2,4).Value = 1.5

$sheetTotal.Cells.Item(3,1).Value = 1
$sheetTotal.Cells.Item(3,2).Value = "2022-Q3"
$sheetTotal.Cells.Item(3,3).Value = 3
$sheetTotal.Cells.Item(3,4).Value = 0.8

$sheetTotal.Cells.Item(4,1).Value = 2
$sheetTotal.Cells.Item(4,2).Value = "2022-Q2"
$sheetTotal.Cells.Item(4,3).Value = 3
$sheetTotal.Cells.Item(4,4).Value = 0.89

$sheetTotal.Cells.Item(2,1).Copy()
$sheetTotal.Cells.Item(5,1).PasteSpecial(-4122)
$sheetTotal.Cells.Item(5,1).Value = 3
$sheetTotal.Cells.Item(5,2).Value = "2022-Q1"
$sheetTotal.Cells.Item(5,3).Value = 1
$sheetTotal.Cells.Item(5,4).Value = 0.67
